$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the existing "keterangan" column (column T / index 20).
# This shifts the old column T (keterangan) to column U, and creates a new blank
# column T for the new "status_kirim_wa" field.
$ws.Columns.Item(20).Insert()

# New column header
$ws.Cells.Item(1, 20).Value = "status_kirim_wa"

# Fill in the new status_kirim_wa values for each shipment row
$ws.Cells.Item(2, 20).Value = "Ya"
$ws.Cells.Item(3, 20).Value = "Ya"
$ws.Cells.Item(4, 20).Value = "Ya"
$ws.Cells.Item(5, 20).Value = "Tidak"
$ws.Cells.Item(6, 20).Value = "Tidak"

# Update the no_resi values (column B) to the new receipt numbers
$ws.Cells.Item(2, 2).Value = "JHD1827183971"
$ws.Cells.Item(3, 2).Value = "JHD1827183972"
$ws.Cells.Item(4, 2).Value = "JHD1827183973"
$ws.Cells.Item(5, 2).Value = "JHD1827183974"
$ws.Cells.Item(6, 2).Value = "JHD1827183975"

# Match the final selection left in the sheet by the author
$ws.Range("C6").Select()
